$d = $word.ActiveDocument

$bullet1 = [char]0x2022 + " Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections"
$bullet2 = [char]0x2022 + " Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government"
$bullet3 = [char]0x2022 + " Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations"

$replacement = "Data Engineering and Infrastructure Architecture^p$bullet1^p$bullet2^p$bullet3"

$found = $d.Content.Find.Execute(
    "Data Engineering and Infrastructure Architecture",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $replacement,
    2
)

Write-Output "Replaced: $found"
